$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new column before G. This shifts the existing "Assignment" (G)
# and "Chart" (H) columns one place to the right (H and I respectively),
# carrying their data/formatting/column-width metadata along with them.
# ---------------------------------------------------------------------------
$ws.Columns("G").Insert()

# New "Data" column width (raw stored width 20 == ColumnWidth 19.1666...
# once run through this host's MDW=6 pixel-grid rounding).
$ws.Columns("G").ColumnWidth = 19.166666666666668

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "Slides"
$ws.Range("F1").Value = "Code"
$ws.Range("G1").Value = "Data"
$ws.Range("H1").Value = "Assignment"
$ws.Range("I1").Value = "Chart"

# ---------------------------------------------------------------------------
# 01_introduction row
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "Introduction"
$ws.Range("E2").Value = "01_introduction"

# ---------------------------------------------------------------------------
# 02_data row
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "Data"
$ws.Range("E3").Value = "02_data"
$ws.Range("F3").Value = "02_data.R"
$ws.Range("G3").Value = "02_data.RData"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = "rstudio"

# ---------------------------------------------------------------------------
# 03_visualization row
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = "Visualization"
$ws.Range("E4").Value = "03_visualization"
$ws.Range("F4").Value = "03_visualization.R"
$ws.Range("G4").Value = "03_visualization.RData"

# ---------------------------------------------------------------------------
# 04_growth row
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Growth · Geometries"
$ws.Range("E5").Value = "04_growth"
$ws.Range("F5").Value = "04_growth.R"
$ws.Range("G5").Value = "04_growth.RData"

# ---------------------------------------------------------------------------
# 05_inflation row
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = "Inflation · Colors"
$ws.Range("E6").Value = "05_inflation"
$ws.Range("F6").Value = "05_inflation.R"
$ws.Range("G6").Value = "05_inflation.RData"

# ---------------------------------------------------------------------------
# 06_labour row (time slot also moved from 10:30-12:30 to 10:00-12:00)
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "10:00-12:00"
$ws.Range("D7").Value = "Labour · Labels"
$ws.Range("E7").Value = "06_labour"
$ws.Range("F7").Value = "06_labour.R"
$ws.Range("G7").Value = "06_labour.RData"
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = "inflation"

# ---------------------------------------------------------------------------
# 07_income row
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "Income · Scales"
$ws.Range("E8").Value = "07_income"
$ws.Range("F8").Value = "07_income.Rmd"
$ws.Range("G8").Value = "07_income.RData"
$ws.Range("G8").Font.Size = 12
$ws.Range("G8").Font.Color = 0

# ---------------------------------------------------------------------------
# Wealth row
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "Wealth · Themes"
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = "incomeineq"

# ---------------------------------------------------------------------------
# Mobility row
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = "Mobility · Maps"

# ---------------------------------------------------------------------------
# Climate row
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = "Climate · Facets"
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = "mobility"

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("B8").Select() | Out-Null
